# Add 2022-Q4 data
# 1) Insert a new summary row in the "总计" (Total) sheet for 2022-Q4.
# 2) Insert a brand-new "2022-Q4" worksheet (before "2022-Q3") holding the
#    per-fund holding detail for that quarter. Inserting it there naturally
#    shifts every later quarter tab (2022-Q3 ... 2020-Q4) one position to
#    the right while keeping their content untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update the "总计" sheet
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Insert a new row right below the header and make sure it doesn't
# inherit the bold header styling that Excel would otherwise copy down.
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

# Give the new index cell (column A) the same style used by the other
# index cells in that column.
$wsTotal.Range("A3").Copy($wsTotal.Range("A2"))

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.57

# Re-number the index column (0,1,2,...) for the rows that shifted down.
for ($r = 3; $r -le 10; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# Step 2: add the new "2022-Q4" sheet with per-fund holdings
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$wsQ4 = $wb.Worksheets.Add($refSheet)
$wsQ4.Name = "2022-Q4"

# Borrow header + row formatting from the existing "2022-Q3" sheet
# (now shifted to position 3) so fonts/borders/column layout match.
$srcSheet = $wb.Worksheets.Item(3)
$srcSheet.Range("B1:H2").Copy($wsQ4.Range("B1"))
$srcSheet.Range("A2").Copy($wsQ4.Range("A2"))
$srcSheet.Range("A2").Copy($wsQ4.Range("A3"))

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# Row 2: fund 000893
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").NumberFormat = "@"
$wsQ4.Range("B2").Value = "000893"
$wsQ4.Range("B2").ClearFormats()
$wsQ4.Range("C2").Value = "工银创新动力股票"
$wsQ4.Range("D2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "15.58"
$wsQ4.Range("D2").ClearFormats()
$wsQ4.Range("E2").NumberFormat = "@"
$wsQ4.Range("E2").Value = "81.84"
$wsQ4.Range("E2").ClearFormats()
$wsQ4.Range("F2").NumberFormat = "@"
$wsQ4.Range("F2").Value = "3.43"
$wsQ4.Range("F2").ClearFormats()
$wsQ4.Range("G2").NumberFormat = "@"
$wsQ4.Range("G2").Value = "0.5344"
$wsQ4.Range("G2").ClearFormats()
$wsQ4.Range("H2").Value = 7

# Row 3: fund 011376
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").NumberFormat = "@"
$wsQ4.Range("B3").Value = "011376"
$wsQ4.Range("B3").ClearFormats()
$wsQ4.Range("C3").Value = "华宝安享混合"
$wsQ4.Range("D3").NumberFormat = "@"
$wsQ4.Range("D3").Value = "4.63"
$wsQ4.Range("D3").ClearFormats()
$wsQ4.Range("E3").NumberFormat = "@"
$wsQ4.Range("E3").Value = "24.21"
$wsQ4.Range("E3").ClearFormats()
$wsQ4.Range("F3").NumberFormat = "@"
$wsQ4.Range("F3").Value = "0.73"
$wsQ4.Range("F3").ClearFormats()
$wsQ4.Range("G3").NumberFormat = "@"
$wsQ4.Range("G3").Value = "0.0338"
$wsQ4.Range("G3").ClearFormats()
$wsQ4.Range("H3").Value = 4

# ---------------------------------------------------------------------
# Keep the original active/selected tab on the "2020-Q4" sheet (it was
# the selected tab before the edit and remains the last sheet now).
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
